# started work on offline caching
# - Log the last week of entries on the "Week 7" sheet (rows 4-10)
# - Add a new "Sheet2" tab after "Week 7" to start this week's log

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Week 7" sheet: fix the B3 date and append the new log rows/columns.
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Week 7")

$ws7.Range("B3").Value = 42807

$ws7.Range("C4").Value = "meeting with Tim, discussed using standard deviation to find outliers in data"

$ws7.Range("B5").Value = 42808
$ws7.Range("B5").NumberFormat = $ws7.Range("B3").NumberFormat
$ws7.Range("C5").Value = "started work on implementing standard deviation "

$ws7.Range("B6").Value = 42809
$ws7.Range("B6").NumberFormat = $ws7.Range("B3").NumberFormat
$ws7.Range("C6").Value = "implemented standard deviation, using the population version instead of sample "

$ws7.Range("B7").Value = 42810
$ws7.Range("B7").NumberFormat = $ws7.Range("B3").NumberFormat
$ws7.Range("C7").Value = "looking into implementing QR code scanning"

$ws7.Range("K4").Value = "issue with graph not displaying all data, cant scroll with dates"

$ws7.Range("B8").Value = 42811
$ws7.Range("B8").NumberFormat = $ws7.Range("B3").NumberFormat
$ws7.Range("C8").Value = "Work on report"

$ws7.Range("B9").Value = 42812
$ws7.Range("B9").NumberFormat = $ws7.Range("B3").NumberFormat
$ws7.Range("C9").Value = "-"

$ws7.Range("B10").Value = 42813
$ws7.Range("B10").NumberFormat = $ws7.Range("B3").NumberFormat
$ws7.Range("C10").Value = "-"

$null = $ws7.Range("K4").Select()

# ---------------------------------------------------------------------
# New "Sheet2" tab after "Week 7" - start of this week's log.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

$ws9.Range("B3").Value = 42814
$ws9.Range("B3").NumberFormat = $ws7.Range("B3").NumberFormat

$null = $ws9.Range("B4").Select()
